# This script applies the updated cryptocurrency market data values
# (price and 1h volume change percentages, plus a few row re-orderings)
# captured in the source XML diff, using Excel COM interop calls only.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '63.673.17'
$ws.Range('E2').Value = '  -0.49%  '

# Row 3
$ws.Range('D3').Value = '2.721.77'
$ws.Range('E3').Value = '  -1.27%  '

# Row 4
$ws.Range('E4').Value = '  -0.06%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '560.35'
$ws.Range('E5').Value = '  -2.53%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '157.96'
$ws.Range('E6').Value = '  -0.68%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.05%  '

# Row 8
$ws.Range('E8').Value = '  -1.52%  '

# Row 9
$ws.Range('E9').Value = '  -2.33%  '

# Row 10
$ws.Range('E10').Value = '  +0.09%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.61'
$ws.Range('E11').Value = '  -3.01%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.374'
$ws.Range('E12').Value = '  -3.19%  '

# Row 13
$ws.Range('D13').Value = '3.199.61'
$ws.Range('E13').Value = '  -1.43%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.60'
$ws.Range('E14').Value = '  -1.37%  '

# Row 15
$ws.Range('D15').Value = '63.523.69'
$ws.Range('E15').Value = '  -0.21%  '

# Row 16
$ws.Range('E16').Value = '  -2.47%  '

# Row 17
$ws.Range('D17').Value = '2.721.23'
$ws.Range('E17').Value = '  -1.52%  '

# Row 18
$ws.Range('E18').Value = '  +0.29%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.66'
$ws.Range('E19').Value = '  -3.91%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '349.93'
$ws.Range('E20').Value = '  -1.77%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.43'
$ws.Range('E21').Value = '  -3.92%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  +0.10%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.515'
$ws.Range('E23').Value = '  -2.45%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.21'
$ws.Range('E24').Value = '  -1.27%  '

# Row 25
$ws.Range('E25').Value = '  +0.28%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.01%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.22'
$ws.Range('E27').Value = '  -4.10%  '

# Row 28
$ws.Range('D28').Value = '0.0₃0885'
$ws.Range('E28').Value = '  -1.58%  '

# Row 29
$ws.Range('E29').Value = '  +10.22%  '

# Row 30
$ws.Range('E30').Value = '  +0.17%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.18'
$ws.Range('E31').Value = '  -2.03%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '164.68'
$ws.Range('E32').Value = '  -2.47%  '

# Row 33
$ws.Range('B33').Value = 'USDe'
$ws.Range('C33').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.998'
$ws.Range('E33').Value = '  -0.04%  '

# Row 34
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '19.87'
$ws.Range('E34').Value = '  -1.16%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.85'
$ws.Range('E35').Value = '  -1.35%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.47'
$ws.Range('E36').Value = '  -0.55%  '

# Row 37
$ws.Range('E37').Value = '  +0.07%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '349.71'
$ws.Range('E38').Value = '  +0.30%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.960'
$ws.Range('E39').Value = '  -4.70%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.15'
$ws.Range('E40').Value = '  -1.51%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.01'
$ws.Range('E41').Value = '  -4.08%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '38.28'
$ws.Range('E42').Value = '  -1.98%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.41'
$ws.Range('E43').Value = '  -1.75%  '

# Row 44
$ws.Range('E44').Value = '  -3.52%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0574'
$ws.Range('E45').Value = '  -2.50%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.626'
$ws.Range('E46').Value = '  -1.07%  '

# Row 47
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '132.40'
$ws.Range('E47').Value = '  -3.28%  '

# Row 48
$ws.Range('B48').Value = 'FirstDigitalUSD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.997'
$ws.Range('E48').Value = '  -0.12%  '

# Row 49
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0246'
$ws.Range('E49').Value = '  -3.11%  '

# Row 50
$ws.Range('B50').Value = 'WhiteBITCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '11.05'
$ws.Range('E50').Value = '  +0.20%  '

# Row 51
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0984'
$ws.Range('E51').Value = '  -2.84%  '
